$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13; this shifts the existing rows 13-30
# down to 14-31 (preserving all of their data/formatting), matching the
# diff which moves every record in that block down by one row and adds
# a brand-new record in the now-empty row 13.
$ws.Rows.Item(13).Insert()

$ws.Range("A13").Value = 1
$ws.Range("B13").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C13").Value = "Arica y Parinacota"
$ws.Range("D13").Value = Get-Date -Year 2023 -Month 3 -Day 9 -Hour 0 -Minute 0 -Second 0
$ws.Range("E13").Value = 15
$ws.Range("F13").Value = "Fruta"
$ws.Range("G13").Value = 100103
$ws.Range("H13").Value = "Frutos de hueso (carozo)"
$ws.Range("I13").Value = 100103002
$ws.Range("J13").Value = "Ciruela"
$ws.Range("K13").Value = "Angeleno"
$ws.Range("L13").Value = "Segunda"
$ws.Range("M13").Value = 270
$ws.Range("N13").Value = 19000
$ws.Range("O13").Value = 20000
$ws.Range("P13").Value = 19500
$ws.Range("Q13").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R13").Value = "Región Metropolitana"
$ws.Range("S13").Value = 1083
$ws.Range("T13").Value = 18

# Apply the same date number format used by the other date cells in column D
$ws.Range("D13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
